# Updates odds values on Sheet1 to match the 2025-02-17 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("G7").Value = 2.4
$ws.Range("I7").Value = 3.6
$ws.Range("L7").Value = 4.5
$ws.Range("N7").Value = 4.75
$ws.Range("AB7").Value = 9.5
$ws.Range("AC7").Value = 12
$ws.Range("AQ7").Value = 51

# Row 8
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("AR8").Value = 4.2
$ws.Range("AS8").Value = 1.23

# Row 11
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 4.2
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 5.5
$ws.Range("M11").Value = 1.14
$ws.Range("N11").Value = 5.5
$ws.Range("U11").Value = 7
$ws.Range("V11").Value = 1.1
$ws.Range("AJ11").Value = 101
$ws.Range("AL11").Value = 7.5
$ws.Range("AM11").Value = 19
$ws.Range("AP11").Value = 41

# Row 18
$ws.Range("M18").Value = 1.11
$ws.Range("N18").Value = 6.5
$ws.Range("O18").Value = 1.57
$ws.Range("P18").Value = 2.25
$ws.Range("S18").Value = 2.7
$ws.Range("T18").Value = 1.44
